# "clean-up of input tables"
# Rename the sheet from "updated" to "Tabelle1" and move the active
# selection from F8 to A5 (matches the captured end-state of the
# workbook after the author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Tabelle1"
$ws.Range("A5").Select()
